$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Remove the last two slides ("Jenkins Pipeline as Code" and
#    "Benefits of Using Jenkins") together with their notes pages.
#    Delete from the end so earlier indices stay stable.
# ---------------------------------------------------------------------------
$p.Slides.Item(10).Delete()
$p.Slides.Item(9).Delete()

# ---------------------------------------------------------------------------
# Each remaining slide keeps its shape order/names (Text 0 = title,
# Text 1.. = bullets); only the text changes, plus two slides (now
# "The Waterfall Model Phases" and "Disadvantages of the Waterfall Model")
# gain extra bullet shapes.
# ---------------------------------------------------------------------------

# Slide 1: Introduction to Jenkins -> Introduction to the Waterfall Model
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Introduction to the Waterfall Model"
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "• A sequential, linear software development approach."
$s1.Shapes.Item(3).TextFrame.TextRange.Text = "• Each phase must be completed before the next can begin."
$s1.Shapes.Item(4).TextFrame.TextRange.Text = "• Focuses on thorough planning and documentation upfront."
$s1.Shapes.Item(5).TextFrame.TextRange.Text = "• Also called a ""Linear Sequential Life Cycle Model""."

# Slide 2: Key Features of Jenkins -> The Waterfall Model Phases (+2 bullets)
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "The Waterfall Model Phases"
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "• **Requirements:** Define project goals and constraints."
$s2.Shapes.Item(3).TextFrame.TextRange.Text = "• **Design:** Plan the architecture and system specifications."
$s2.Shapes.Item(4).TextFrame.TextRange.Text = "• **Implementation:** Code and build the system based on the design."
$s2.Shapes.Item(5).TextFrame.TextRange.Text = "• **Testing:** Verify and validate the system meets requirements."

$s2t5 = $s2.Shapes.AddTextbox(1, 36, 216, 540, 0)
$s2t5.Name = "Text 5"
$s2t5.Height = 0
$s2t5.Fill.Visible = $false
$s2t5tf = $s2t5.TextFrame
$s2t5tf.WordWrap = $true
$s2t5tf.VerticalAnchor = 3
$s2t5tr = $s2t5tf.TextRange
$s2t5tr.Text = "• **Deployment:** Release the system to the users."
$s2t5tr.Font.Size = 20
$s2t5tr.Font.Color.RGB = 0
$s2t5tr.ParagraphFormat.Bullet.Visible = $false

$s2t6 = $s2.Shapes.AddTextbox(1, 36, 252, 540, 0)
$s2t6.Name = "Text 6"
$s2t6.Height = 0
$s2t6.Fill.Visible = $false
$s2t6tf = $s2t6.TextFrame
$s2t6tf.WordWrap = $true
$s2t6tf.VerticalAnchor = 3
$s2t6tr = $s2t6tf.TextRange
$s2t6tr.Text = "• **Maintenance:** Fix bugs and improve the system."
$s2t6tr.Font.Size = 20
$s2t6tr.Font.Color.RGB = 0
$s2t6tr.ParagraphFormat.Bullet.Visible = $false

# Slide 3: Jenkins Architecture -> Requirements Phase
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Requirements Phase"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "• A detailed definition of system functionality."
$s3.Shapes.Item(3).TextFrame.TextRange.Text = "• Output is a comprehensive requirements document."
$s3.Shapes.Item(4).TextFrame.TextRange.Text = "• Key questions: What is the system supposed to do? What problem does it solve?"
$s3.Shapes.Item(5).TextFrame.TextRange.Text = "• User stories can play an important part here."

# Slide 4: Continuous Integration with Jenkins -> Design Phase
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Design Phase"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "• Transforms requirements into a detailed design."
$s4.Shapes.Item(3).TextFrame.TextRange.Text = "• Defines system architecture, data structures, and algorithms."
$s4.Shapes.Item(4).TextFrame.TextRange.Text = "• Outputs include design documents and specifications."
$s4.Shapes.Item(5).TextFrame.TextRange.Text = "• Covers both high-level and low-level design details."

# Slide 5: Continuous Delivery with Jenkins -> Implementation Phase
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Implementation Phase"
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "• Translates the design into executable code."
$s5.Shapes.Item(3).TextFrame.TextRange.Text = "• Coding is the primary activity in this phase."
$s5.Shapes.Item(4).TextFrame.TextRange.Text = "• Unit testing is often performed during implementation."
$s5.Shapes.Item(5).TextFrame.TextRange.Text = "• Requires adherence to coding standards."

# Slide 6: Setting Up Jenkins -> Testing Phase
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Testing Phase"
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "• Evaluates the system against the requirements."
$s6.Shapes.Item(3).TextFrame.TextRange.Text = "• Detects and fixes defects in the software."
$s6.Shapes.Item(4).TextFrame.TextRange.Text = "• Various testing methods employed (unit, integration, system)."
$s6.Shapes.Item(5).TextFrame.TextRange.Text = "• Produces test reports and bug tracking information."

# Slide 7: Creating a Jenkins Job -> Advantages of the Waterfall Model
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Advantages of the Waterfall Model"
$s7.Shapes.Item(2).TextFrame.TextRange.Text = "• Simple and easy to understand and implement."
$s7.Shapes.Item(3).TextFrame.TextRange.Text = "• Well-suited for projects with clear, stable requirements."
$s7.Shapes.Item(4).TextFrame.TextRange.Text = "• Disciplined approach with defined stages and milestones."
$s7.Shapes.Item(5).TextFrame.TextRange.Text = "• Easy to manage due to the rigidity of the model."

# Slide 8: Jenkins Plugins -> Disadvantages of the Waterfall Model (+1 bullet)
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Disadvantages of the Waterfall Model"
$s8.Shapes.Item(2).TextFrame.TextRange.Text = "• Inflexible to changes once a phase is complete."
$s8.Shapes.Item(3).TextFrame.TextRange.Text = "• Difficult to accommodate evolving requirements."
$s8.Shapes.Item(4).TextFrame.TextRange.Text = "• Working software is only produced late in the lifecycle."
$s8.Shapes.Item(5).TextFrame.TextRange.Text = "• High risk and uncertainty if requirements are not fully understood upfront."

$s8t5 = $s8.Shapes.AddTextbox(1, 36, 216, 540, 0)
$s8t5.Name = "Text 5"
$s8t5.Height = 0
$s8t5.Fill.Visible = $false
$s8t5tf = $s8t5.TextFrame
$s8t5tf.WordWrap = $true
$s8t5tf.VerticalAnchor = 3
$s8t5tr = $s8t5tf.TextRange
$s8t5tr.Text = "• Not a good model for complex or object-oriented projects."
$s8t5tr.Font.Size = 20
$s8t5tr.Font.Color.RGB = 0
$s8t5tr.ParagraphFormat.Bullet.Visible = $false
